$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = '["01", "04", "08", "10", "17", "22", "26", "30"]'
$ws.Range("B20").Value = "18:15 - 18:19"
$ws.Range("B21").Value = "18:20 - 18:24"
$ws.Range("C21").Value = '["01", "04", "08", "10", "17", "22", "26", "30"]'

[void]$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
